{"js": "// Thesis ch.06 (\"Simulace\" section on the Metropolis-Hastings algorithm) edits:\n//\n// 1) \", jinak \" + \"\\u03b8\" used to be two separate <w:r> runs with identical\n//    (default) formatting; Word re-merges adjacent same-formatted runs like\n//    this whenever the text around them is edited, so normalize them back\n//    into a single run.\n// 2) Likewise \"...hustoto\" + \"u pravd\\u011bpodobnosti...\" used to be split\n//    into two runs around a stray leftover \"_GoBack\" bookmark; merge that\n//    text back into a single run too.\n// 3) Italicize the term \"random walk\" (it is the English gloss for\n//    \"n\\u00e1hodn\\u00e1 ch\\u016fze\" and should render in italics like the other\n//    foreign/variable terms in this thesis).\n// 4) + 5) Word only ever keeps one \"_GoBack\" bookmark (the last place the\n//    user edited). Move it from its old spot (inside the \"hustotou\"\n//    sentence) to the new edit location, right before \"dn\\u00e9 ch\\u016fze (\" in\n//    the random-walk sentence.\n\n// --- 1) merge \", jinak \" + \"\\u03b8\" -------------------------------------------\n{\n  const target = \", jinak \\u03b8\";\n  const results = context.document.body.search(target, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  results.items[0].insertText(target, \"Replace\");\n  await context.sync();\n}\n\n// --- 2) merge \"...hustoto\" + \"u pravd\\u011bpodobnosti...\" --------------------\n{\n  const target =\n    \"Toto vede k\\u00a0tomu, \\u017ee oblasti s\\u00a0vy\\u0161\\u0161\\u00ed hustotou \" +\n    \"pravd\\u011bpodobnosti jsou n\\u00e1hodnou ch\\u016fz\\u00ed nav\\u0161t\\u00edveny \" +\n    \"v\\u00edcekr\\u00e1t, a proto, pokud vybereme dostate\\u010dn\\u00e9 mno\\u017estv\\u00ed \" +\n    \"vzork\\u016f (a zbav\\u00edme se \";\n  const results = context.document.body.search(target, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  results.items[0].insertText(target, \"Replace\");\n  await context.sync();\n}\n\n// --- 3) italicize \"random walk\" ------------------------------------------\n{\n  const results = context.document.body.search(\"random walk\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  results.items[0].font.italic = true;\n  await context.sync();\n}\n\n// --- 4) drop whatever is left of the old \"_GoBack\" bookmark --------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 5) re-create \"_GoBack\" right before \"dn\\u00e9 ch\\u016fze (\" -------------------\n{\n  const anchor = \"pomoc\\u00ed n\\u00e1ho\";\n  const results = context.document.body.search(anchor, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const insertionPoint = results.items[0].getRange(\"End\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Thesis ch.06 (\"Simulace\" section on the Metropolis-Hastings algorithm) edits:\n#\n# 1) \", jinak \" + the following theta used to be two separate runs with\n#    identical (default) formatting; Word re-merges adjacent same-formatted\n#    runs like this whenever the text around them is edited, so normalize\n#    them back into a single run.\n# 2) Likewise \"...hustoto\" + \"u pravdepodobnosti...\" used to be split into\n#    two runs around a stray leftover \"_GoBack\" bookmark; merge that text\n#    back into a single run too.\n# 3) Italicize the term \"random walk\" (it is the English gloss for the\n#    Czech \"nahodna chuze\" and should render in italics like the other\n#    foreign/variable terms in this thesis).\n# 4) + 5) Word only ever keeps one \"_GoBack\" bookmark (the last place the\n#    user edited). Move it from its old spot (inside the \"hustotou\"\n#    sentence) to the new edit location, right before \"dne chuze (\" in the\n#    random-walk sentence.\n#\n# NOTE: the Czech strings below are built from explicit character codes\n# ([char]0x....) instead of being typed literally, to guarantee the exact\n# characters (incl. the U+00A0 NBSPs) survive byte-for-byte, chunked to\n# stay under the interpreter's expression-nesting limit.\n\n$d = $word.ActiveDocument\n\n# --- 1) merge \", jinak \" + theta ------------------------------------------\n$mergeJinak = (\"\" + [char]0x002C + [char]0x0020 + [char]0x006A + [char]0x0069 + [char]0x006E + [char]0x0061 + [char]0x006B + [char]0x0020 + [char]0x03B8)\n$range1 = $d.Content\n$range1.Find.Execute($mergeJinak, $false, $false, $false, $false, $false, $true, 1, $false, $mergeJinak, 2) | Out-Null\n\n# --- 2) merge \"...hustoto\" + \"u pravdepodobnosti...\" ----------------------\n$mergeHustotu = (\"\" + [char]0x0054 + [char]0x006F + [char]0x0074 + [char]0x006F + [char]0x0020 + [char]0x0076 + [char]0x0065 + [char]0x0064 + [char]0x0065 + [char]0x0020 + [char]0x006B + [char]0x00A0 + [char]0x0074 + [char]0x006F + [char]0x006D + [char]0x0075 + [char]0x002C + [char]0x0020 + [char]0x017E + [char]0x0065 + [char]0x0020 + [char]0x006F + [char]0x0062 + [char]0x006C + [char]0x0061 + [char]0x0073 + [char]0x0074 + [char]0x0069 + [char]0x0020 + [char]0x0073 + [char]0x00A0 + [char]0x0076 + [char]0x0079 + [char]0x0161 + [char]0x0161 + [char]0x00ED + [char]0x0020 + [char]0x0068 + [char]0x0075 + [char]0x0073 + [char]0x0074 + [char]0x006F + [char]0x0074 + [char]0x006F + [char]0x0075 + [char]0x0020 + [char]0x0070 + [char]0x0072 + [char]0x0061 + [char]0x0076 + [char]0x0064 + [char]0x011B + [char]0x0070 + [char]0x006F + [char]0x0064 + [char]0x006F + [char]0x0062 + [char]0x006E + [char]0x006F + [char]0x0073 + [char]0x0074 + [char]0x0069 + [char]0x0020 + [char]0x006A + [char]0x0073 + [char]0x006F + [char]0x0075 + [char]0x0020 + [char]0x006E + [char]0x00E1 + [char]0x0068 + [char]0x006F + [char]0x0064 + [char]0x006E + [char]0x006F + [char]0x0075 + [char]0x0020 + [char]0x0063 + [char]0x0068 + [char]0x016F) + (\"\" + [char]0x007A + [char]0x00ED + [char]0x0020 + [char]0x006E + [char]0x0061 + [char]0x0076 + [char]0x0161 + [char]0x0074 + [char]0x00ED + [char]0x0076 + [char]0x0065 + [char]0x006E + [char]0x0079 + [char]0x0020 + [char]0x0076 + [char]0x00ED + [char]0x0063 + [char]0x0065 + [char]0x006B + [char]0x0072 + [char]0x00E1 + [char]0x0074 + [char]0x002C + [char]0x0020 + [char]0x0061 + [char]0x0020 + [char]0x0070 + [char]0x0072 + [char]0x006F + [char]0x0074 + [char]0x006F + [char]0x002C + [char]0x0020 + [char]0x0070 + [char]0x006F + [char]0x006B + [char]0x0075 + [char]0x0064 + [char]0x0020 + [char]0x0076 + [char]0x0079 + [char]0x0062 + [char]0x0065 + [char]0x0072 + [char]0x0065 + [char]0x006D + [char]0x0065 + [char]0x0020 + [char]0x0064 + [char]0x006F + [char]0x0073 + [char]0x0074 + [char]0x0061 + [char]0x0074 + [char]0x0065 + [char]0x010D + [char]0x006E + [char]0x00E9 + [char]0x0020 + [char]0x006D + [char]0x006E + [char]0x006F + [char]0x017E + [char]0x0073 + [char]0x0074 + [char]0x0076 + [char]0x00ED + [char]0x0020 + [char]0x0076 + [char]0x007A + [char]0x006F + [char]0x0072 + [char]0x006B + [char]0x016F + [char]0x0020 + [char]0x0028 + [char]0x0061 + [char]0x0020 + [char]0x007A + [char]0x0062) + (\"\" + [char]0x0061 + [char]0x0076 + [char]0x00ED + [char]0x006D + [char]0x0065 + [char]0x0020 + [char]0x0073 + [char]0x0065 + [char]0x0020)\n$range2 = $d.Content\n$range2.Find.Execute($mergeHustotu, $false, $false, $false, $false, $false, $true, 1, $false, $mergeHustotu, 2) | Out-Null\n\n# --- 3) italicize \"random walk\" -------------------------------------------\n$range3 = $d.Content\n$range3.Find.Execute(\"random walk\") | Out-Null\n$range3.Italic = 1\n\n# --- 4) drop whatever is left of the old \"_GoBack\" bookmark ---------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 5) re-create \"_GoBack\" right before \"dne chuze (\" --------------------\n$anchor = (\"\" + [char]0x0070 + [char]0x006F + [char]0x006D + [char]0x006F + [char]0x0063 + [char]0x00ED + [char]0x0020 + [char]0x006E + [char]0x00E1 + [char]0x0068 + [char]0x006F)\n$range5 = $d.Content\n$range5.Find.Execute($anchor) | Out-Null\n$goBackRange = $range5.Duplicate\n$goBackRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n"}
